$p = $ppt.ActivePresentation

# --- Slide 1: "Rectangle 27" -- "Justen and Liam" -> "All" ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "Person doing task: Justen and Liam") {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Length
            $tail = $tr.Characters($len - 14, 15)
            $tail.Text = "All"
        }
    }
}

# --- Slide 2: "Rectangle 31" -- ":" -> ": Justen" ---
# --- Slide 2: "Rectangle 33" -- split "Person doing task" / ":" -> ": Liam" ---
$s2 = $p.Slides.Item(2)
$seen = 0
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "Person doing task:") {
            $seen = $seen + 1
            $tr = $shp.TextFrame.TextRange
            if ($seen -eq 1) {
                # First occurrence -> append "Justen" after the colon.
                $colon = $tr.Characters(18, 1)
                $colon.Text = ": Justen"
            } else {
                # Second occurrence -> split "task" into its own run, then
                # append "Liam" after the colon.
                $taskRun = $tr.Characters(14, 4)
                $taskRun.Text = "task"
                $colon = $tr.Characters(18, 1)
                $colon.Text = ": Liam"
            }
        }
    }
}
